$d = $word.ActiveDocument

$d.Content.Find.Execute("19×47=", $true, $false, $false, $false, $false, $true, 1, $false, "37×47=", 2) | Out-Null
$d.Content.Find.Execute("57×52=", $true, $false, $false, $false, $false, $true, 1, $false, "69×37=", 2) | Out-Null
$d.Content.Find.Execute("25×51=", $true, $false, $false, $false, $false, $true, 1, $false, "46×29=", 2) | Out-Null
$d.Content.Find.Execute("55×95=", $true, $false, $false, $false, $false, $true, 1, $false, "76×30=", 2) | Out-Null
$d.Content.Find.Execute("24×66=", $true, $false, $false, $false, $false, $true, 1, $false, "71×81=", 2) | Out-Null
$d.Content.Find.Execute("28×87=", $true, $false, $false, $false, $false, $true, 1, $false, "46×65=", 2) | Out-Null
$d.Content.Find.Execute("96×98=", $true, $false, $false, $false, $false, $true, 1, $false, "58×77=", 2) | Out-Null
$d.Content.Find.Execute("59×76=", $true, $false, $false, $false, $false, $true, 1, $false, "95×22=", 2) | Out-Null
$d.Content.Find.Execute("92×57=", $true, $false, $false, $false, $false, $true, 1, $false, "35×41=", 2) | Out-Null
$d.Content.Find.Execute("69×58=", $true, $false, $false, $false, $false, $true, 1, $false, "90×78=", 2) | Out-Null
$d.Content.Find.Execute("95×71=", $true, $false, $false, $false, $false, $true, 1, $false, "12×94=", 2) | Out-Null
$d.Content.Find.Execute("24×61=", $true, $false, $false, $false, $false, $true, 1, $false, "99×76=", 2) | Out-Null
$d.Content.Find.Execute("35×79=", $true, $false, $false, $false, $false, $true, 1, $false, "26×34=", 2) | Out-Null
$d.Content.Find.Execute("69×18=", $true, $false, $false, $false, $false, $true, 1, $false, "53×56=", 2) | Out-Null
$d.Content.Find.Execute("18×67=", $true, $false, $false, $false, $false, $true, 1, $false, "26×73=", 2) | Out-Null
$d.Content.Find.Execute("23×16=", $true, $false, $false, $false, $false, $true, 1, $false, "39×85=", 2) | Out-Null
$d.Content.Find.Execute("93×42=", $true, $false, $false, $false, $false, $true, 1, $false, "38×26=", 2) | Out-Null
$d.Content.Find.Execute("38×30=", $true, $false, $false, $false, $false, $true, 1, $false, "15×47=", 2) | Out-Null
$d.Content.Find.Execute("35×85=", $true, $false, $false, $false, $false, $true, 1, $false, "72×93=", 2) | Out-Null
$d.Content.Find.Execute("36×91=", $true, $false, $false, $false, $false, $true, 1, $false, "14×35=", 2) | Out-Null
$d.Content.Find.Execute("63×17=", $true, $false, $false, $false, $false, $true, 1, $false, "75×41=", 2) | Out-Null
$d.Content.Find.Execute("34×64=", $true, $false, $false, $false, $false, $true, 1, $false, "76×53=", 2) | Out-Null
$d.Content.Find.Execute("83×61=", $true, $false, $false, $false, $false, $true, 1, $false, "45×50=", 2) | Out-Null
$d.Content.Find.Execute("92×72=", $true, $false, $false, $false, $false, $true, 1, $false, "40×55=", 2) | Out-Null
$d.Content.Find.Execute("81×26=", $true, $false, $false, $false, $false, $true, 1, $false, "13×98=", 2) | Out-Null
